# Add a new "2022" column (column S) to the right of the existing "2021"
# column (column R), copying formatting from column R and filling in the
# 2022 values for each region row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> 2022 value (row 4 is the header year row)
$values = @{
    4  = 2022
    5  = 4.9538761752705343
    6  = 11.304954640614097
    7  = 5.1593323216995444
    8  = 13.687943262411348
    9  = 10.22864019253911
    10 = 9.1213700670141478
    11 = 3.1335149863760217
    12 = 2.872905173311127
    13 = 3.527842284697861
    14 = 5.0305321314335565
}

foreach ($row in 4..14) {
    $src = $ws.Range("R$row")
    $dst = $ws.Range("S$row")

    # Copy the formatting (number format, font, borders, alignment) from the
    # corresponding 2021 (column R) cell onto the new 2022 (column S) cell.
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats

    # Write in the actual 2022 figure.
    $dst.Value = $values[$row]
}

$excel.CutCopyMode = 0

# Match the recorded selection state after the edit.
$ws.Range("T6").Select() | Out-Null
